# Auto-generated edit script applying numeric corrections to market-price derived columns
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2216.6667
$ws.Range("J17").Value = 2216.6667
$ws.Range("L17").Value = 6650.000100000001
$ws.Range("N17").Value = -6986.000100000001

# Row 40
$ws.Range("H40").Value = 5350.4546
$ws.Range("I40").Value = 3499
$ws.Range("J40").Value = 5535.6
$ws.Range("K40").Value = 3499
$ws.Range("L40").Value = 5535.6
$ws.Range("M40").Value = -3324
$ws.Range("N40").Value = -5885.6

# Row 55
$ws.Range("H55").Value = 581.7857
$ws.Range("I55").Value = 558.7
$ws.Range("K55").Value = 558.7
$ws.Range("M55").Value = -344.7

# Row 74
$ws.Range("H74").Value = 5940.75
$ws.Range("I74").Value = 4754.3335
$ws.Range("K74").Value = 4754.3335
$ws.Range("M74").Value = -3818.3335

# Row 77
$ws.Range("H77").Value = 5940.75
$ws.Range("I77").Value = 4754.3335
$ws.Range("K77").Value = 23771.6675
$ws.Range("M77").Value = -19091.6675

# Row 80
$ws.Range("H80").Value = 4360.905
$ws.Range("I80").Value = 4540.143
$ws.Range("J80").Value = 4271.2856
$ws.Range("K80").Value = 13620.429
$ws.Range("L80").Value = 12813.8568
$ws.Range("M80").Value = -12622.429
$ws.Range("N80").Value = -14809.8568

# Row 83
$ws.Range("H83").Value = 4360.905
$ws.Range("I83").Value = 4540.143
$ws.Range("J83").Value = 4271.2856
$ws.Range("K83").Value = 40861.287
$ws.Range("L83").Value = 38441.5704
$ws.Range("M83").Value = -35869.287
$ws.Range("N83").Value = -48425.5704

# Row 95
$ws.Range("H95").Value = 20884.666
$ws.Range("J95").Value = 20884.666
$ws.Range("L95").Value = 20884.666
$ws.Range("N95").Value = -26376.666

# Row 98
$ws.Range("H98").Value = 2670.3635
$ws.Range("I98").Value = 2737.5
$ws.Range("K98").Value = 2737.5
$ws.Range("M98").Value = -1239.5

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

# Row 118
$ws.Range("H118").Value = 720.7143
$ws.Range("I118").Value = 720.7143
$ws.Range("K118").Value = 2162.1429
$ws.Range("M118").Value = -505.1428999999998

# Row 122
$ws.Range("H122").Value = 2670.3635
$ws.Range("I122").Value = 2737.5
$ws.Range("K122").Value = 8212.5
$ws.Range("M122").Value = -5762.5

# Row 138
$ws.Range("H138").Value = 3420.861
$ws.Range("I138").Value = 1657.2941
$ws.Range("J138").Value = 4998.7896
$ws.Range("K138").Value = 4971.8823
$ws.Range("L138").Value = 14996.3688
$ws.Range("M138").Value = 168.1176999999998
$ws.Range("N138").Value = -25276.3688

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2829.5
$ws.Range("I74").Value = 3055.6667
$ws.Range("K74").Value = 3055.6667
$ws.Range("M74").Value = -2181.6667

# Row 77
$ws.Range("H77").Value = 2829.5
$ws.Range("I77").Value = 3055.6667
$ws.Range("K77").Value = 15278.3335
$ws.Range("M77").Value = -10910.3335

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 167.83333
$ws.Range("I22").Value = 161.33333
$ws.Range("K22").Value = 161.33333
$ws.Range("M22").Value = 11.66667000000001

# Row 80
$ws.Range("H80").Value = 5062.769
$ws.Range("J80").Value = 9143
$ws.Range("L80").Value = 9143
$ws.Range("N80").Value = -11139

# Row 83
$ws.Range("H83").Value = 5062.769
$ws.Range("J83").Value = 9143
$ws.Range("L83").Value = 45715
$ws.Range("N83").Value = -55699

# Row 86
$ws.Range("H86").Value = 4666.5
$ws.Range("I86").Value = 4750
$ws.Range("J86").Value = 4499.5
$ws.Range("K86").Value = 4750
$ws.Range("L86").Value = 4499.5
$ws.Range("M86").Value = -3627
$ws.Range("N86").Value = -6745.5

# Row 89
$ws.Range("H89").Value = 4666.5
$ws.Range("I89").Value = 4750
$ws.Range("J89").Value = 4499.5
$ws.Range("K89").Value = 23750
$ws.Range("L89").Value = 22497.5
$ws.Range("M89").Value = -18134
$ws.Range("N89").Value = -33729.5

# Row 105
$ws.Range("H105").Value = 3475
$ws.Range("I105").Value = 2674.5715
$ws.Range("K105").Value = 2674.5715
$ws.Range("M105").Value = -927.5715

# Row 107
$ws.Range("H107").Value = 4333
$ws.Range("I107").Value = 5000
$ws.Range("J107").Value = 3999.5
$ws.Range("K107").Value = 5000
$ws.Range("L107").Value = 3999.5
$ws.Range("M107").Value = -3080
$ws.Range("N107").Value = -7839.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1280.6666
$ws.Range("I16").Value = 1246
$ws.Range("J16").Value = 1350
$ws.Range("K16").Value = 1246
$ws.Range("L16").Value = 1350
$ws.Range("M16").Value = -959
$ws.Range("N16").Value = -1924

# Row 60
$ws.Range("H60").Value = 36356.855
$ws.Range("I60").Value = 18625
$ws.Range("K60").Value = 18625
$ws.Range("M60").Value = -18114

# Row 113
$ws.Range("H113").Value = 1280.6666
$ws.Range("I113").Value = 1246
$ws.Range("J113").Value = 1350
$ws.Range("K113").Value = 1246
$ws.Range("L113").Value = 1350
$ws.Range("M113").Value = 924
$ws.Range("N113").Value = -5690

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 335665.5
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# Row 39
$ws.Range("H39").Value = 28500
$ws.Range("J39").Value = 28500
$ws.Range("L39").Value = 28500
$ws.Range("N39").Value = -29564

# Row 96
$ws.Range("H96").Value = 39999
$ws.Range("J96").Value = 39999
$ws.Range("L96").Value = 39999
$ws.Range("N96").Value = -45491

# Row 105
$ws.Range("H105").Value = 22443.166
$ws.Range("J105").Value = 22443.166
$ws.Range("L105").Value = 22443.166
$ws.Range("N105").Value = -29431.166

# Row 107
$ws.Range("H107").Value = 2988.5715
$ws.Range("I107").Value = 1184.2
$ws.Range("K107").Value = 1184.2
$ws.Range("M107").Value = 735.8

# Row 113
$ws.Range("H113").Value = 1666.3334
$ws.Range("I113").Value = 1666.3334
$ws.Range("K113").Value = 1666.3334
$ws.Range("M113").Value = 503.6666

# Row 132
$ws.Range("H132").Value = 3634.1667
$ws.Range("I132").Value = 4161
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 12483
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -9953
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 16620.525
$ws.Range("I40").Value = 14523.929
$ws.Range("K40").Value = 14523.929
$ws.Range("M40").Value = -14387.929

# Row 46
$ws.Range("H46").Value = 2203.3103
$ws.Range("I46").Value = 1831
$ws.Range("K46").Value = 1831
$ws.Range("M46").Value = -1643

# Row 68
$ws.Range("H68").Value = 2500
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751

# Row 71
$ws.Range("H71").Value = 2500
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756

# Row 87
$ws.Range("H87").Value = 26250
$ws.Range("J87").Value = 26250
$ws.Range("L87").Value = 26250
$ws.Range("N87").Value = -28496

# Row 90
$ws.Range("H90").Value = 26250
$ws.Range("J90").Value = 26250
$ws.Range("L90").Value = 78750
$ws.Range("N90").Value = -89982

# Row 122
$ws.Range("H122").Value = 15155.429
$ws.Range("I122").Value = 14247.9
$ws.Range("J122").Value = 17424.25
$ws.Range("K122").Value = 42743.7
$ws.Range("L122").Value = 52272.75
$ws.Range("M122").Value = -40293.7
$ws.Range("N122").Value = -57172.75

# Row 136
$ws.Range("H136").Value = 1411.5
$ws.Range("I136").Value = 1248.6666
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 3745.9998
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -1195.9998
$ws.Range("N136").Value = -10800

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 6238.25
$ws.Range("I132").Value = 4414.4287
$ws.Range("K132").Value = 13243.2861
$ws.Range("M132").Value = -10713.2861

# Row 136
$ws.Range("H136").Value = 8913.362999999999
$ws.Range("I136").Value = 8167.75
$ws.Range("J136").Value = 10901.667
$ws.Range("K136").Value = 24503.25
$ws.Range("L136").Value = 32705.001
$ws.Range("M136").Value = -21953.25
$ws.Range("N136").Value = -37805.001
